$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet right after the existing one
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Copy the data table (header + 31 days) from sheet1 into the new sheet, preserving formatting
$src = $ws1.Range("A9:K40")
$src.Copy($ws2.Range("A1"))

# Update selection on the source sheet to reflect the copied range
[void]$ws1.Range("A9:K40").Select()

# Update selection on the new sheet and make it the active tab
[void]$ws2.Range("A1:K32").Select()
$ws2.Activate()

Write-Output "done"
